$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 110, shifting existing rows 110-209 down to 111-210.
$ws.Rows.Item(110).Insert()

# Populate the newly inserted row 110 with the new weekly price-report record.
$ws.Cells.Item(110, 1).Value = 4
$ws.Cells.Item(110, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(110, 3).Value = "Los Lagos"
$ws.Cells.Item(110, 4).Value = 44566
$ws.Cells.Item(110, 5).Value = 10
$ws.Cells.Item(110, 6).Value = 100112040
$ws.Cells.Item(110, 7).Value = "Cilantro"
$ws.Cells.Item(110, 8).Value = "Sin especificar"
$ws.Cells.Item(110, 9).Value = "Primera"
$ws.Cells.Item(110, 10).Value = 20
$ws.Cells.Item(110, 11).Value = 12000
$ws.Cells.Item(110, 12).Value = 12000
$ws.Cells.Item(110, 13).Value = 12000
$ws.Cells.Item(110, 14).Value = "$/docena de atados (2 kilos)"
$ws.Cells.Item(110, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(110, 16).Value = 6000
$ws.Cells.Item(110, 17).Value = 2
$ws.Cells.Item(110, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(110, 4).NumberFormat = $ws.Cells.Item(111, 4).NumberFormat
